# issue #5: stock data from json to db
# The "股票" (stock) sheet gains a "category" column (value "normal" for every
# row) right after "property_category", and two trailing columns
# "source_file" ("tmpa5201" for every row) and "index" (a copy of the
# original row index already stored in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 12

# 1) Insert a new column before the existing "date" column (I) and turn it
#    into the new "category" column. This shifts date/legislator_name/
#    legislator_id one column to the right (I->J, J->K, K->L) while keeping
#    their formatting intact.
$ws.Columns("I:I").Insert()

$ws.Range("I1").Value = "category"
$ws.Range("I2:I$lastRow").Value = "normal"

# 2) Append the two brand new trailing columns: "source_file" and "index".
#    Copy the header/data styles from the neighbouring "legislator_id"
#    column (now L) so the new cells keep matching formatting.
$ws.Range("L1").Copy($ws.Range("M1"))
$ws.Range("L1").Copy($ws.Range("N1"))
$ws.Range("L2:L$lastRow").Copy($ws.Range("M2:M$lastRow"))
$ws.Range("L2:L$lastRow").Copy($ws.Range("N2:N$lastRow"))

$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$ws.Range("M2:M$lastRow").Value = "tmpa5201"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value
}
